$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Edmund"
$ws.Range("B2").Value = "January"
$ws.Range("C2").Value = "Brave Crystal-Dazzler"

$ws.Range("B2").Select() | Out-Null
